# Updated cryptos list on Tue Aug 20 16:13:29 UTC 2024 with GitHub Actions
# Refreshes Price (D) / Volume(1h) (E) figures, and re-ranks the
# Aave / EnergySwap pair (rows 48-49) to match the new scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values are plain decimals (e.g. "1.00", "6.00") that Excel
# would otherwise silently coerce to numbers (dropping the trailing zero)
# when assigned through .Value. Force those specific cells to Text format
# first so they stay literal strings, matching the source data.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '58.723.73'
$ws.Range("E2").Value = '  +0.44%  '
$ws.Range("D3").Value = '2.567.58'
$ws.Range("E3").Value = '  -0.44%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '559.39'
$ws.Range("E5").Value = '  +2.85%  '
$ws.Range("D6").Value = '142.08'
$ws.Range("E6").Value = '  -1.17%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  +2.13%  '
$ws.Range("D9").Value = '2.573.33'
$ws.Range("E9").Value = '  -0.53%  '
$ws.Range("E10").Value = '  -1.83%  '
$ws.Range("E11").Value = '  +2.36%  '
$ws.Range("E12").Value = '  +8.38%  '
$ws.Range("D13").Value = '0.340'
$ws.Range("E13").Value = '  +2.04%  '
$ws.Range("D14").Value = '3.021.56'
$ws.Range("E14").Value = '  -0.55%  '
$ws.Range("D15").Value = '58.804.46'
$ws.Range("E15").Value = '  +0.69%  '
$ws.Range("E16").Value = '  +5.92%  '
$ws.Range("E17").Value = '  +3.68%  '
$ws.Range("D18").Value = '2.573.48'
$ws.Range("E18").Value = '  -0.62%  '
$ws.Range("D19").Value = '4.48'
$ws.Range("E19").Value = '  +0.48%  '
$ws.Range("D20").Value = '333.97'
$ws.Range("D21").Value = '10.11'
$ws.Range("E21").Value = '  +0.75%  '
$ws.Range("D22").Value = '6.15'
$ws.Range("E22").Value = '  +1.20%  '
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("D24").Value = '63.57'
$ws.Range("E24").Value = '  -4.08%  '
$ws.Range("E25").Value = '  +4.75%  '
$ws.Range("E26").Value = '  +0.24%  '
$ws.Range("D27").Value = '0.161'
$ws.Range("E27").Value = '  +1.78%  '
$ws.Range("E28").Value = '  +1.53%  '
$ws.Range("D29").Value = '0.0₃0775'
$ws.Range("E29").Value = '  +5.16%  '
$ws.Range("D30").Value = '0.998'
$ws.Range("E30").Value = '  -0.01%  '
$ws.Range("D31").Value = '1.66'
$ws.Range("E31").Value = '  +1.12%  '
$ws.Range("D32").Value = '160.15'
$ws.Range("E32").Value = '  +4.11%  '
$ws.Range("D33").Value = '6.00'
$ws.Range("E33").Value = '  +0.67%  '
$ws.Range("E34").Value = '  -0.13%  '
$ws.Range("E35").Value = '  +2.07%  '
$ws.Range("D36").Value = '0.877'
$ws.Range("E36").Value = '  +2.86%  '
$ws.Range("D37").Value = '0.871'
$ws.Range("E37").Value = '  +6.40%  '
$ws.Range("E38").Value = '  +2.39%  '
$ws.Range("E39").Value = '  -0.92%  '
$ws.Range("E40").Value = '  +3.53%  '
$ws.Range("D41").Value = '291.70'
$ws.Range("E41").Value = '  +4.49%  '
$ws.Range("D42").Value = '3.60'
$ws.Range("E42").Value = '  +0.84%  '
$ws.Range("D43").Value = '1.00'
$ws.Range("E44").Value = '  +3.48%  '
$ws.Range("D45").Value = '0.590'
$ws.Range("E45").Value = '  -0.56%  '
$ws.Range("D46").Value = '10.60'
$ws.Range("E46").Value = '  -0.35%  '
$ws.Range("D47").Value = '0.0531'
$ws.Range("E47").Value = '  +0.60%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").Value = '124.66'
$ws.Range("E48").Value = '  +12.39%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '18.86'
$ws.Range("E49").Value = '  +1.84%  '
$ws.Range("E50").Value = '  +0.98%  '
$ws.Range("D51").Value = '18.34'
$ws.Range("E51").Value = '  +3.22%  '
